# Fix profile metric value is 0 — the N column (time metric) was being
# computed at half its correct value; double N22:N43 to correct it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 22; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 14)  # column N is the 14th column
    $cell.Value = $cell.Value2 * 2
}
